# "Running temps and updated log"
#
# - Row 49/50 B column: the old ad-hoc timestamp names
#   ("2022-06-23 16-03-44" / "2022-06-23 16-04-56") are replaced by the
#   pipeline run name "Pipe_SCTv2_23-06" (matching rows 42-48).
# - 4 new "post selection" DEG rows are appended (53-56).
# - 4 new SingleR Rdata rows are appended (57-60) documenting the newest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- update existing rows 49 & 50 (col B only) ---
$ws.Range("B49").Value = "Pipe_SCTv2_23-06"
$ws.Range("B50").Value = "Pipe_SCTv2_23-06"

# --- new DEG "post selection" rows ---
$ws.Range("C53").Value = "DEG"
$ws.Range("D53").Value = "A + C"
$ws.Range("E53").Value = "old postSelection"

$ws.Range("C54").Value = "DEG"
$ws.Range("D54").Value = "A + C"
$ws.Range("E54").Value = "new postSelection"

$ws.Range("C55").Value = "DEG"
$ws.Range("D55").Value = "N + C"
$ws.Range("E55").Value = "old postSelection"

$ws.Range("C56").Value = "DEG"
$ws.Range("D56").Value = "N + C"
$ws.Range("E56").Value = "new postSelection"

# --- new SingleR Rdata rows for the 2022-06-24 run ---
$ws.Range("A57").Value = "Rdata "
$ws.Range("C57").Value = "Kriegstein to SingleR"
$ws.Range("D57").Value = "A+C"
$ws.Range("E57").Value = "oldPostSelect"
$ws.Range("B57").Value = "SingleR_RData_2022-06-24 08-55-17"

$ws.Range("A58").Value = "Rdata "
$ws.Range("C58").Value = "Kriegstein to SingleR"
$ws.Range("D58").Value = "A+C"
$ws.Range("E58").Value = "oldSelect"
$ws.Range("B58").Value = "SingleR_RData_2022-06-24 08-58-57"

$ws.Range("A59").Value = "Rdata "
$ws.Range("C59").Value = "Kriegstein to SingleR"
$ws.Range("D59").Value = "A+C"
$ws.Range("E59").Value = "newPostSelect"

$ws.Range("A60").Value = "Rdata "
$ws.Range("C60").Value = "Kriegstein to SingleR"
$ws.Range("D60").Value = "A+C"
$ws.Range("E60").Value = "newSelect"

$ws.Range("B59").Value = "SingleR_RData_2022-06-24 09-04-40"
$ws.Range("B60").Value = "SingleR_RData_2022-06-24 09-07-27"

# --- move the active selection to reflect where editing finished ---
$ws.Activate()
$ws.Range("C65").Select()
